$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (TestCase A / Headers): add header*X-User=jack to the header string ---
$ws.Range("G3").Value = "header*hdr1=ABC&header*hdr2=DEF&header*X-User=jack"

# --- Row 4 (TestCase A / Expected): add X-EntryPoint and X-User to expected json ---
$ws.Range("G4").Value = '{"User":"jack","name":"jack","role":"admin","Host":"localhost","hdr1":"ABC","hdr2":"DEF","X-EntryPoint":"TestProject","X-HostPath":"localhost","X-User":"jack"}'

# --- Row 7 (TestCase B / Expected): add X-EntryPoint to expected json ---
$ws.Range("G7").Value = '{"User":"jill","role":"user","group":"456","Host":"localhost","hdr1":"123","X-EntryPoint":"TestProject","X-User":"jill"}'

# --- Row 8 (TestCase C / Claims): claim prefix added to the claim string ---
$ws.Range("G8").Value = "claim*X-User=bob"

# --- Row 10 (TestCase C / Expected): expected json now matches jill instead of bob ---
$ws.Range("G10").Value = '{"User":"jill"}'

# --- Widen column G to fit the newly lengthened content ---
$ws.Columns.Item(7).ColumnWidth = 46.75

# --- Move the active selection from G10 to G9 ---
$ws.Range("G9").Select() | Out-Null
